$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "Method"
$ws.Range("B5").Value = 0.9949715634984005
$ws.Range("C5").Value = 0.9949967221633021
$ws.Range("D5").Value = 147340769.6732198
$ws.Range("E5").Value = 5509.050019388131
$ws.Range("F5").Value = 0.9634307661821588
$ws.Range("G5").Value = 0.9985709573263347
$ws.Range("H5").Value = 1.361060513652995
$ws.Range("I5").Value = 1.361060513652995
$ws.Range("J5").Value = 1.390264645806614
$ws.Range("K5").Value = 106822
